$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3_")

# Update existing A2/B2 values
$ws.Range("A2").Value = 0.93
$ws.Range("B2").Value = 0.03

# Enter new header/label strings in the same order the original author
# typed them, so the shared-strings table is built in matching order.
$ws.Range("F1").Value = "reading"
$ws.Range("G1").Value = "Fitted value"
$ws.Range("J1").Value = "K"
$ws.Range("E1").Value = "Mass"
$ws.Range("H1").Value = "Deviation"
$ws.Range("L1").Value = "DoF"
$ws.Range("K1").Value = "N"

# Data rows
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.55
$ws.Range("J2").Value = 0.55
$ws.Range("G2").Formula = "=E2*`$J`$2"
$ws.Range("H2").Formula = "=(F2-G2)^2"
$ws.Range("K2").Value = 4
$ws.Range("L2").Formula = "=K2-2"

$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1.06
$ws.Range("G3").Formula = "=E3*`$J`$2"
$ws.Range("H3").Formula = "=(F3-G3)^2"

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1.71
$ws.Range("G4").Formula = "=E4*`$J`$2"
$ws.Range("H4").Formula = "=(F4-G4)^2"

$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 2.09
$ws.Range("G5").Formula = "=E5*`$J`$2"
$ws.Range("H5").Formula = "=(F5-G5)^2"

$ws.Range("G8").Value = "SoF"
$ws.Range("H8").Formula = "=SQRT(SUM(H2:H5)/L2)"

# Update sheet view: set active cell selection on this sheet
$ws.Range("B3").Select()

# Make this the active sheet (tab selected)
$ws.Activate()
